$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Cells.Item(1, 6).Value = "time_taken"

$ws.Cells.Item(2, 6).Value = "2021-10-05 13:40:21.271690"
$ws.Cells.Item(3, 6).Value = "2021-10-05 13:40:21.271702"
$ws.Cells.Item(4, 6).Value = "2021-10-05 13:40:21.271706"
$ws.Cells.Item(5, 6).Value = "2021-10-05 13:40:21.271710"
$ws.Cells.Item(6, 6).Value = "2021-10-05 13:40:21.271713"
$ws.Cells.Item(7, 6).Value = "2021-10-05 13:40:21.271716"
$ws.Cells.Item(8, 6).Value = "2021-10-05 13:40:21.271720"
$ws.Cells.Item(9, 6).Value = "2021-10-05 13:40:21.271723"
$ws.Cells.Item(10, 6).Value = "2021-10-05 13:40:21.271726"
$ws.Cells.Item(11, 6).Value = "2021-10-05 13:40:21.271729"
$ws.Cells.Item(12, 6).Value = "2021-10-05 13:40:21.271732"
$ws.Cells.Item(13, 6).Value = "2021-10-05 13:40:21.271760"
$ws.Cells.Item(14, 6).Value = "2021-10-05 13:40:21.271768"
$ws.Cells.Item(15, 6).Value = "2021-10-05 13:40:21.271771"
$ws.Cells.Item(16, 6).Value = "2021-10-05 13:40:21.271774"
$ws.Cells.Item(17, 6).Value = "2021-10-05 13:40:21.271778"
$ws.Cells.Item(18, 6).Value = "2021-10-05 13:40:21.271782"
$ws.Cells.Item(19, 6).Value = "2021-10-05 13:40:21.271785"
$ws.Cells.Item(20, 6).Value = "2021-10-05 13:40:21.271788"
$ws.Cells.Item(21, 6).Value = "2021-10-05 13:40:21.271791"
$ws.Cells.Item(22, 6).Value = "2021-10-05 13:40:21.271795"
$ws.Cells.Item(23, 6).Value = "2021-10-05 13:40:21.271798"
$ws.Cells.Item(24, 6).Value = "2021-10-05 13:40:21.271801"
$ws.Cells.Item(25, 6).Value = "2021-10-05 13:40:21.271804"
$ws.Cells.Item(26, 6).Value = "2021-10-05 13:40:21.271808"
$ws.Cells.Item(27, 6).Value = "2021-10-05 13:40:21.271811"
$ws.Cells.Item(28, 6).Value = "2021-10-05 13:40:21.271814"
$ws.Cells.Item(29, 6).Value = "2021-10-05 13:40:21.271817"
$ws.Cells.Item(30, 6).Value = "2021-10-05 13:40:21.271820"
$ws.Cells.Item(31, 6).Value = "2021-10-05 13:40:21.271823"
$ws.Cells.Item(32, 6).Value = "2021-10-05 13:40:21.271826"
$ws.Cells.Item(33, 6).Value = "2021-10-05 13:40:21.271829"
$ws.Cells.Item(34, 6).Value = "2021-10-05 13:40:21.271833"
$ws.Cells.Item(35, 6).Value = "2021-10-05 13:40:21.271836"
$ws.Cells.Item(36, 6).Value = "2021-10-05 13:40:21.271839"
$ws.Cells.Item(37, 6).Value = "2021-10-05 13:40:21.271842"
$ws.Cells.Item(38, 6).Value = "2021-10-05 13:40:21.271845"
$ws.Cells.Item(39, 6).Value = "2021-10-05 13:40:21.271848"
$ws.Cells.Item(40, 6).Value = "2021-10-05 13:40:21.271852"
$ws.Cells.Item(41, 6).Value = "2021-10-05 13:40:21.271855"
$ws.Cells.Item(42, 6).Value = "2021-10-05 13:40:21.271858"
$ws.Cells.Item(43, 6).Value = "2021-10-05 13:40:21.271862"
$ws.Cells.Item(44, 6).Value = "2021-10-05 13:40:21.271865"
$ws.Cells.Item(45, 6).Value = "2021-10-05 13:40:21.271868"
